# Auto-generated: apply Universalis market-data refresh to all 8 Leve-profit sheets.
# For each changed cell we assign the new literal value; two cells whose column
# ('M') became inapplicable after the refresh are cleared instead (matches the diff,
# which drops those <c> elements entirely rather than writing 0).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 846.0769
$ws.Range("I2").Value = 799.9
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 799.9
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -686.9
$ws.Range("N2").Value = -1226
$ws.Range("H62").Value = 2886.7896
$ws.Range("I62").Value = 2331.2727
$ws.Range("J62").Value = 3650.625
$ws.Range("K62").Value = 2331.2727
$ws.Range("L62").Value = 3650.625
$ws.Range("M62").Value = -1707.2727
$ws.Range("N62").Value = -4898.625
$ws.Range("H65").Value = 2886.7896
$ws.Range("I65").Value = 2331.2727
$ws.Range("J65").Value = 3650.625
$ws.Range("K65").Value = 11656.3635
$ws.Range("L65").Value = 18253.125
$ws.Range("M65").Value = -8536.363499999999
$ws.Range("N65").Value = -24493.125
$ws.Range("H76").Value = 4633.3335
$ws.Range("I76").Value = 3156.5715
$ws.Range("J76").Value = 6700.8
$ws.Range("K76").Value = 3156.5715
$ws.Range("L76").Value = 6700.8
$ws.Range("M76").Value = -2841.5715
$ws.Range("N76").Value = -7330.8
$ws.Range("H79").Value = 4633.3335
$ws.Range("I79").Value = 3156.5715
$ws.Range("J79").Value = 6700.8
$ws.Range("K79").Value = 3156.5715
$ws.Range("L79").Value = 6700.8
$ws.Range("M79").Value = -2064.5715
$ws.Range("N79").Value = -8884.799999999999
$ws.Range("H100").Value = 66667984
$ws.Range("I100").Value = 66667984
$ws.Range("K100").Value = 66667984
$ws.Range("M100").Value = -66667443
$ws.Range("H131").Value = 4142
$ws.Range("I131").Value = 4229
$ws.Range("J131").Value = 4098.5
$ws.Range("K131").Value = 12687
$ws.Range("L131").Value = 12295.5
$ws.Range("M131").Value = -7647
$ws.Range("N131").Value = -22375.5
$ws.Range("H138").Value = 2584.98
$ws.Range("I138").Value = 849.4286
$ws.Range("J138").Value = 2867.5117
$ws.Range("K138").Value = 2548.2858
$ws.Range("L138").Value = 8602.535100000001
$ws.Range("M138").Value = 2591.7142
$ws.Range("N138").Value = -18882.5351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 23538.25
$ws.Range("J24").Value = 23538.25
$ws.Range("L24").Value = 23538.25
$ws.Range("N24").Value = -24286.25
$ws.Range("H32").Value = 4206.915
$ws.Range("I32").Value = 3476.9092
$ws.Range("J32").Value = 14244.5
$ws.Range("K32").Value = 3476.9092
$ws.Range("L32").Value = 14244.5
$ws.Range("M32").Value = -3189.9092
$ws.Range("N32").Value = -14818.5
$ws.Range("H63").Value = 17317356
$ws.Range("I63").Value = 34630212
$ws.Range("J63").Value = 4499.25
$ws.Range("K63").Value = 34630212
$ws.Range("L63").Value = 4499.25
$ws.Range("M63").Value = -34629526
$ws.Range("N63").Value = -5871.25
$ws.Range("H66").Value = 17317356
$ws.Range("I66").Value = 34630212
$ws.Range("J66").Value = 4499.25
$ws.Range("K66").Value = 173151060
$ws.Range("L66").Value = 22496.25
$ws.Range("M66").Value = -173147628
$ws.Range("N66").Value = -29360.25
$ws.Range("H88").Value = 11113794
$ws.Range("J88").Value = 3266.6667
$ws.Range("L88").Value = 3266.6667
$ws.Range("N88").Value = -4078.6667
$ws.Range("H91").Value = 11113794
$ws.Range("J91").Value = 3266.6667
$ws.Range("L91").Value = 3266.6667
$ws.Range("N91").Value = -6074.6667
$ws.Range("H100").Value = 23538.25
$ws.Range("J100").Value = 23538.25
$ws.Range("L100").Value = 23538.25
$ws.Range("N100").Value = -25702.25
$ws.Range("H132").Value = 1766.4412
$ws.Range("I132").Value = 880.92
$ws.Range("K132").Value = 2642.76
$ws.Range("M132").Value = -112.7599999999998
$ws.Range("H133").Value = 19531.8
$ws.Range("J133").Value = 19531.8
$ws.Range("L133").Value = 19531.8
$ws.Range("N133").Value = -24591.8
$ws.Range("H137").Value = 41030
$ws.Range("J137").Value = 41030
$ws.Range("L137").Value = 41030
$ws.Range("N137").Value = -51230
$ws.Range("H139").Value = 40920.77
$ws.Range("J139").Value = 40920.77
$ws.Range("L139").Value = 40920.77
$ws.Range("N139").Value = -51200.77

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 23666.334
$ws.Range("J44").Value = 23666.334
$ws.Range("L44").Value = 23666.334
$ws.Range("N44").Value = -24660.334
$ws.Range("H86").Value = 2999.8
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2999.8
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2999.8
$ws.Range("N86").Value = -5245.8
$ws.Range("H89").Value = 2999.8
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2999.8
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 14999
$ws.Range("N89").Value = -26231
$ws.Range("H137").Value = 45730
$ws.Range("J137").Value = 45730
$ws.Range("L137").Value = 45730
$ws.Range("N137").Value = -55930
$ws.Range("M86").ClearContents()
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6899748.5
$ws.Range("I99").Value = 15386161
$ws.Range("K99").Value = 15386161
$ws.Range("M99").Value = -15384663
$ws.Range("H126").Value = 6899748.5
$ws.Range("I126").Value = 15386161
$ws.Range("K126").Value = 46158483
$ws.Range("M126").Value = -46156013
$ws.Range("H134").Value = 1638.6774
$ws.Range("I134").Value = 913.6818
$ws.Range("K134").Value = 2741.0454
$ws.Range("M134").Value = -206.0454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4464895
$ws.Range("J113").Value = 9615969
$ws.Range("L113").Value = 28847907
$ws.Range("N113").Value = -28852247

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 27890
$ws.Range("J42").Value = 27890
$ws.Range("L42").Value = 27890
$ws.Range("N42").Value = -28860
$ws.Range("H46").Value = 35151.2
$ws.Range("J46").Value = 35151.2
$ws.Range("L46").Value = 35151.2
$ws.Range("N46").Value = -35463.2
$ws.Range("H63").Value = 15250
$ws.Range("J63").Value = 15250
$ws.Range("L63").Value = 15250
$ws.Range("N63").Value = -16622
$ws.Range("H66").Value = 15250
$ws.Range("J66").Value = 15250
$ws.Range("L66").Value = 45750
$ws.Range("N66").Value = -52614
$ws.Range("H70").Value = 6348.5312
$ws.Range("I70").Value = 5680.4
$ws.Range("J70").Value = 7462.0835
$ws.Range("K70").Value = 5680.4
$ws.Range("L70").Value = 7462.0835
$ws.Range("M70").Value = -5410.4
$ws.Range("N70").Value = -8002.0835
$ws.Range("H73").Value = 6348.5312
$ws.Range("I73").Value = 5680.4
$ws.Range("J73").Value = 7462.0835
$ws.Range("K73").Value = 5680.4
$ws.Range("L73").Value = 7462.0835
$ws.Range("M73").Value = -4744.4
$ws.Range("N73").Value = -9334.083500000001
$ws.Range("H80").Value = 50002890
$ws.Range("I80").Value = 62502860
$ws.Range("K80").Value = 62502860
$ws.Range("M80").Value = -62501862
$ws.Range("H83").Value = 50002890
$ws.Range("I83").Value = 62502860
$ws.Range("K83").Value = 312514300
$ws.Range("M83").Value = -312509308
$ws.Range("H113").Value = 1729.2142
$ws.Range("I113").Value = 1708.3846
$ws.Range("K113").Value = 1708.3846
$ws.Range("M113").Value = 461.6153999999999
$ws.Range("H115").Value = 27890
$ws.Range("J115").Value = 27890
$ws.Range("L115").Value = 27890
$ws.Range("N115").Value = -30240
$ws.Range("H122").Value = 4871.7
$ws.Range("I122").Value = 1859.8
$ws.Range("J122").Value = 7883.6
$ws.Range("K122").Value = 5579.4
$ws.Range("L122").Value = 23650.8
$ws.Range("M122").Value = -3129.4
$ws.Range("N122").Value = -28550.8
$ws.Range("H137").Value = 37125
$ws.Range("J137").Value = 37125
$ws.Range("L137").Value = 37125
$ws.Range("N137").Value = -47325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3357.8
$ws.Range("I7").Value = 2465.3125
$ws.Range("J7").Value = 4944.4443
$ws.Range("K7").Value = 2465.3125
$ws.Range("L7").Value = 4944.4443
$ws.Range("M7").Value = -2353.3125
$ws.Range("N7").Value = -5168.4443
$ws.Range("H40").Value = 7475.6924
$ws.Range("I40").Value = 6214
$ws.Range("J40").Value = 8557.143
$ws.Range("K40").Value = 6214
$ws.Range("L40").Value = 8557.143
$ws.Range("M40").Value = -6078
$ws.Range("N40").Value = -8829.143
$ws.Range("H126").Value = 3357.8
$ws.Range("I126").Value = 2465.3125
$ws.Range("J126").Value = 4944.4443
$ws.Range("K126").Value = 7395.9375
$ws.Range("L126").Value = 14833.3329
$ws.Range("M126").Value = -4925.9375
$ws.Range("N126").Value = -19773.3329
$ws.Range("H133").Value = 28970
$ws.Range("J133").Value = 28970
$ws.Range("L133").Value = 28970
$ws.Range("N133").Value = -34030

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 26900
$ws.Range("J64").Value = 26900
$ws.Range("L64").Value = 26900
$ws.Range("N64").Value = -27396
$ws.Range("H67").Value = 26900
$ws.Range("J67").Value = 26900
$ws.Range("L67").Value = 26900
$ws.Range("N67").Value = -28616
$ws.Range("H133").Value = 43999
$ws.Range("J133").Value = 43999
$ws.Range("L133").Value = 43999
$ws.Range("N133").Value = -54119
